# Generate Report for handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) values for the
# c44b5bcc-... row (row 4) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D4").Value = "2016-01-28 10:52:34"
$zhcn.Range("G4").Value = "2016-01-28 10:53:22"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D4").Value = "2016-01-28 10:52:47"
$dede.Range("G4").Value = "2016-01-28 10:53:44"
